# The draft contained a block of "inspiration quotes" / placeholder
# sections ("[recent works]", "[intro design]", "[design]",
# "[concluzii]" and the paragraphs that went with them) that followed
# the paragraph ending in "... elevi." and ran through to the very end
# of the document ("... scenarii de confuzie."). Remove that whole
# block, leaving the "... elevi." paragraph as the last paragraph of
# the document body.

$d = $word.ActiveDocument

# Paragraph to keep as the new last paragraph: the one ending in
# "... elevi.".
$startFind = $d.Content
$startFind.Find.Execute("elevi.") | Out-Null
$lastKeptPara = $startFind.Paragraphs(1)

# Last paragraph of the document (the one ending in "... scenarii de
# confuzie."), i.e. the end of the block to remove.
$endFind = $d.Content
$endFind.Find.Execute("scenarii de confuzie.") | Out-Null
$lastRemovedPara = $endFind.Paragraphs(1)

# Range spanning from just after the paragraph mark that ends the
# "elevi." paragraph through to (and including) the paragraph mark
# that ends the document's final paragraph.
$toDelete = $d.Range($lastKeptPara.Range.End, $lastRemovedPara.Range.End)
$toDelete.Delete()
